# Drop e Create atualizados
# Adds three new "drop table" statements (livros_colecao, autor_livro,
# livro_idioma) in two places:
#   1. Right before the first "drop table pais;" line.
#   2. Right after "drop table progressoleitura;" (in the second block).

$d = $word.ActiveDocument

$newLines = "drop table livros_colecao;`r`ndrop table autor_livro;`r`ndrop table livro_idioma;`r`n"

# --- Insertion 1: before "drop table pais;" ---
$r1 = $d.Content
$r1.Find.Execute("drop table pais;", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertRange1 = $d.Range($r1.Start, $r1.Start)
$insertRange1.InsertBefore($newLines)

# --- Insertion 2: after "drop table progressoleitura;" ---
$r2 = $d.Content
$r2.Find.Execute("drop table progressoleitura;", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint2 = $r2.End + 1
$insertRange2 = $d.Range($insertPoint2, $insertPoint2)
$insertRange2.InsertBefore($newLines)
